$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, 109.1786205742166, 10.37097476852108, 1.233906563344711),
    @(1, 112.511197023974, 10.23179269695487, 1.088284220444473),
    @(2, 102.7756044855822, 10.58245397251616, 1.054717717983383),
    @(3, 103.166746438103, 10.8996482801407, 0.7926528662174265),
    @(4, 104.6677834151466, 9.923241123419512, 1.213748783451909),
    @(5, 106.7794747921097, 9.487946435860174, 1.323166101186889),
    @(6, 98.54296266253425, 10.69053567380967, 1.064762899394597),
    @(7, 104.3718034604765, 10.33323711615574, 0.9435102513415738),
    @(8, 105.502578037155, 9.693702029651408, 1.223471251108252),
    @(9, 101.1011844393747, 10.37965696862676, 1.066588561399274),
    @(10, 97.03337002844215, 10.76193278594143, 1.006258967614246),
    @(11, 97.64265315771384, 10.29382548146205, 1.124525108184251),
    @(12, 99.46655366285279, 10.52890109413067, 0.8283571727685592),
    @(13, 106.4945079987856, 9.330237392450396, 1.124623496196694),
    @(14, 93.15156922950189, 10.58962240117504, 0.858600470404917),
    @(15, 92.78555524585381, 10.13239105666403, 1.09887750072958),
    @(16, 116.9548527908364, 8.749768031377881, 0.5059689150226316),
    @(17, 88.81350597853542, 10.06948390580832, 1.109296649182347),
    @(18, 88.86949133669478, 9.774395448161119, 1.227896890778952),
    @(19, 90.42192385422325, 10.15493461409919, 0.9148563577195951),
    @(20, 89.38019988050237, 10.45565062792593, 0.755169316139251),
    @(21, 98.42976962175062, 8.51550648349269, 1.24404509069062),
    @(22, 84.82381565628384, 9.797461041534991, 1.230507208966279),
    @(23, 102.6259428699427, 9.289074331226985, 0.4979650740737364),
    @(24, 84.6488430913671, 10.11986427361115, 0.9047655866026827),
    @(25, 92.38981654337677, 8.878086946979682, 1.13158097321209),
    @(26, 82.1327797047794, 9.948318190070859, 1.048045808153214),
    @(27, 81.35088578706497, 10.21020303845059, 0.9203340417559835),
    @(28, 82.89554911232271, 9.712510102896687, 0.6370699908969363),
    @(29, 89.54688633585391, 8.514859368099582, 0.7736691004446931),
    @(30, 84.08595203661967, 8.134077062256477, 1.020412268517304),
    @(31, 90.87304998666903, 7.902647772436168, 0.7552897245731691)
)

for ($i = 0; $i -lt 32; $i++) {
    $row = $data[$i]
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

Write-Host "Done writing data"
